# Ata de Sprint 04-10-2023 - add sentence about Thursday meeting after the
# "Product Backlog" bullet paragraph.
$d = $word.ActiveDocument

# Locate the end of the existing sentence that ends the target paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("Backlog de acordo com a aula de Tecnologia da Informação.", `
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if (-not $found) {
    throw "Could not find anchor text for the Sprint paragraph edit."
}

# Collapse the found range to its end so subsequent inserts land right after
# the existing final period, inside the same paragraph.
$rng.Collapse(0)

# 1) Plain run: " Na "
$rng.InsertAfter(" Na ")
$rng.Collapse(0)

# 2) Bold run: "quinta,"
$boldStart = $rng.End
$rng.InsertAfter("quinta,")
$boldEnd = $rng.End
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Font.Bold = $true
$boldRange.Font.BoldBi = $true
$rng.Collapse(0)

# 3) Plain run: remainder of the sentence
$rng.InsertAfter(" os membros devem se reunir para definir as prioridades e dificuldades de cada requisito.")
